$wb = $excel.ActiveWorkbook

# Sheet 1: Number_Features_All
$ws1 = $wb.Worksheets.Item("Number_Features_All")
$ws1.Range("B2").Value = 65
$ws1.Range("D2").Value = 455

# Sheet 2: Number_Features_Cleaned
$ws2 = $wb.Worksheets.Item("Number_Features_Cleaned")
$ws2.Range("B2").Value = 62
$ws2.Range("D2").Value = 427
